$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 1379.8572
$ws.Range("I122").Value = 1763.5555
$ws.Range("J122").Value = 1198.1052
$ws.Range("K122").Value = 5290.666499999999
$ws.Range("L122").Value = 3594.3156
$ws.Range("M122").Value = -2840.666499999999
$ws.Range("N122").Value = -8494.3156
$ws.Range("H132").Value = 2541.6584
$ws.Range("I132").Value = 1583.1428
$ws.Range("J132").Value = 8133
$ws.Range("K132").Value = 4749.428400000001
$ws.Range("L132").Value = 24399
$ws.Range("M132").Value = -2219.428400000001
$ws.Range("N132").Value = -29459

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 1347.2916
$ws.Range("I134").Value = 1340.579
$ws.Range("J134").Value = 1372.8
$ws.Range("K134").Value = 4021.737
$ws.Range("L134").Value = 4118.4
$ws.Range("M134").Value = -1486.737
$ws.Range("N134").Value = -9188.4

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 49.166668
$ws.Range("I7").Value = 35.25
$ws.Range("J7").Value = 77
$ws.Range("K7").Value = 35.25
$ws.Range("L7").Value = 77
$ws.Range("M7").Value = 77.75
$ws.Range("N7").Value = -303
$ws.Range("H17").Value = 21666.666
$ws.Range("I17").Value = 5000
$ws.Range("J17").Value = 30000
$ws.Range("K17").Value = 5000
$ws.Range("L17").Value = 30000
$ws.Range("M17").Value = -4826
$ws.Range("N17").Value = -30348
$ws.Range("H25").Value = 3265.8
$ws.Range("I25").Value = 587
$ws.Range("J25").Value = 7284
$ws.Range("K25").Value = 587
$ws.Range("L25").Value = 7284
$ws.Range("M25").Value = -413
$ws.Range("N25").Value = -7632
$ws.Range("H31").Value = 2006.3036
$ws.Range("I31").Value = 1397.7778
$ws.Range("J31").Value = 3101.65
$ws.Range("K31").Value = 1397.7778
$ws.Range("L31").Value = 3101.65
$ws.Range("M31").Value = -1102.7778
$ws.Range("N31").Value = -3691.65
$ws.Range("H34").Value = 2006.3036
$ws.Range("I34").Value = 1397.7778
$ws.Range("J34").Value = 3101.65
$ws.Range("K34").Value = 1397.7778
$ws.Range("L34").Value = 3101.65
$ws.Range("M34").Value = -1195.7778
$ws.Range("N34").Value = -3505.65
$ws.Range("H51").Value = 8337
$ws.Range("J51").Value = 9386.4
$ws.Range("L51").Value = 9386.4
$ws.Range("N51").Value = -10858.4
$ws.Range("H59").Value = 14450.833
$ws.Range("I59").Value = 3104
$ws.Range("J59").Value = 16720.2
$ws.Range("K59").Value = 3104
$ws.Range("L59").Value = 16720.2
$ws.Range("M59").Value = -1959
$ws.Range("N59").Value = -19010.2
$ws.Range("H60").Value = 22063.555
$ws.Range("J60").Value = 24434.875
$ws.Range("L60").Value = 24434.875
$ws.Range("N60").Value = -25456.875
$ws.Range("H61").Value = 8337
$ws.Range("J61").Value = 9386.4
$ws.Range("L61").Value = 9386.4
$ws.Range("N61").Value = -10082.4
$ws.Range("H68").Value = 17400
$ws.Range("J68").Value = 17400
$ws.Range("L68").Value = 17400
$ws.Range("N68").Value = -18898
$ws.Range("H71").Value = 17400
$ws.Range("J71").Value = 17400
$ws.Range("L71").Value = 52200
$ws.Range("N71").Value = -59688
$ws.Range("H74").Value = 13773.272
$ws.Range("J74").Value = 13773.272
$ws.Range("L74").Value = 13773.272
$ws.Range("N74").Value = -15521.272
$ws.Range("H77").Value = 13773.272
$ws.Range("J77").Value = 13773.272
$ws.Range("L77").Value = 41319.81600000001
$ws.Range("N77").Value = -50055.81600000001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H57").Value = 20327.5
$ws.Range("J57").Value = 35600
$ws.Range("L57").Value = 35600
$ws.Range("N57").Value = -37240
$ws.Range("H97").Value = 2466.4666
$ws.Range("I97").Value = 2867.5
$ws.Range("J97").Value = 1664.4
$ws.Range("K97").Value = 2867.5
$ws.Range("L97").Value = 1664.4
$ws.Range("M97").Value = -2371.5
$ws.Range("N97").Value = -2656.4
$ws.Range("H102").Value = 2634.4119
$ws.Range("I102").Value = 1539.1818
$ws.Range("J102").Value = 4642.3335
$ws.Range("K102").Value = 1539.1818
$ws.Range("L102").Value = 4642.3335
$ws.Range("M102").Value = 82.81819999999993
$ws.Range("N102").Value = -7886.3335
$ws.Range("H122").Value = 1625.2222
$ws.Range("I122").Value = 1703.375
$ws.Range("J122").Value = 1000
$ws.Range("K122").Value = 5110.125
$ws.Range("L122").Value = 3000
$ws.Range("M122").Value = -2660.125
$ws.Range("N122").Value = -7900
$ws.Range("H126").Value = 12143.538
$ws.Range("I126").Value = 2627.4546
$ws.Range("J126").Value = 19122
$ws.Range("K126").Value = 7882.3638
$ws.Range("L126").Value = 57366
$ws.Range("M126").Value = -5412.3638
$ws.Range("N126").Value = -62306
$ws.Range("H135").Value = 48812.273
$ws.Range("J135").Value = 48812.273
$ws.Range("L135").Value = 48812.273
$ws.Range("N135").Value = -58952.273

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2729.8
$ws.Range("I7").Value = 2412.25
$ws.Range("J7").Value = 4000
$ws.Range("K7").Value = 2412.25
$ws.Range("L7").Value = 4000
$ws.Range("M7").Value = -2300.25
$ws.Range("N7").Value = -4224
$ws.Range("H46").Value = 1158.421
$ws.Range("I46").Value = 1070
$ws.Range("J46").Value = 1238
$ws.Range("K46").Value = 1070
$ws.Range("L46").Value = 1238
$ws.Range("M46").Value = -882
$ws.Range("N46").Value = -1614
$ws.Range("H93").Value = 2519.7273
$ws.Range("I93").Value = 1970.7778
$ws.Range("K93").Value = 1970.7778
$ws.Range("M93").Value = -722.7778000000001
$ws.Range("H122").Value = 3282.5
$ws.Range("I122").Value = 2499.1428
$ws.Range("J122").Value = 4065.8572
$ws.Range("K122").Value = 7497.428400000001
$ws.Range("L122").Value = 12197.5716
$ws.Range("M122").Value = -5047.428400000001
$ws.Range("N122").Value = -17097.5716
$ws.Range("H126").Value = 2729.8
$ws.Range("I126").Value = 2412.25
$ws.Range("J126").Value = 4000
$ws.Range("K126").Value = 7236.75
$ws.Range("L126").Value = 12000
$ws.Range("M126").Value = -4766.75
$ws.Range("N126").Value = -16940
$ws.Range("H132").Value = 3548.7742
$ws.Range("I132").Value = 3514.7727
$ws.Range("J132").Value = 3631.889
$ws.Range("K132").Value = 10544.3181
$ws.Range("L132").Value = 10895.667
$ws.Range("M132").Value = -8014.3181
$ws.Range("N132").Value = -15955.667

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 7925
$ws.Range("I62").Value = 3778.5715
$ws.Range("J62").Value = 17600
$ws.Range("K62").Value = 3778.5715
$ws.Range("L62").Value = 17600
$ws.Range("M62").Value = -3154.5715
$ws.Range("N62").Value = -18848
$ws.Range("H65").Value = 7925
$ws.Range("I65").Value = 3778.5715
$ws.Range("J65").Value = 17600
$ws.Range("K65").Value = 18892.8575
$ws.Range("L65").Value = 88000
$ws.Range("M65").Value = -15772.8575
$ws.Range("N65").Value = -94240
$ws.Range("H75").Value = 29793.8
$ws.Range("J75").Value = 29793.8
$ws.Range("L75").Value = 29793.8
$ws.Range("N75").Value = -31665.8
$ws.Range("H78").Value = 29793.8
$ws.Range("J78").Value = 29793.8
$ws.Range("L78").Value = 89381.39999999999
$ws.Range("N78").Value = -98741.39999999999
$ws.Range("H109").Value = 32600
$ws.Range("J109").Value = 32600
$ws.Range("L109").Value = 32600
$ws.Range("N109").Value = -35374

Write-Host "Applied all changes"